# Puerto Rico_Converted.xlsx - "Updated policies and graphs"
#
# Changes applied:
#  1. "Temp Policy_Industries" (column R) weight goes from 1 to 0 (row 7),
#     which drops the total weight (AC7) from 13 to 12.
#  2. Because AC7 changed, every LockdownEffectiveness value (column AC,
#     rows 20-221) that depends on it is recomputed as
#         AC{row} = SUMPRODUCT(B{row}:AB{row}, $B$7:$AB$7) / $AC$7
#  3. Twelve new daily rows are appended (rows 222-233) for 9/30/2020
#     through 10/11/2020, carrying forward the same closure pattern as
#     the last existing row (221) and the corresponding recomputed
#     LockdownEffectiveness value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- 1. Zero out the "Temp Policy_Industries" weight and the new total ---
$ws.Range("R7").Value = 0
$ws.Range("AC7").Value = 12

# --- 2. Recalculated LockdownEffectiveness values for existing rows ---
$ws.Range("AC20:AC23").Value = 0.08333333333333333
$ws.Range("AC24:AC72").Value = 0.758333333325
$ws.Range("AC73:AC94").Value = 0.8416666666583333
$ws.Range("AC95:AC95").Value = 0.8138888888833332
$ws.Range("AC96:AC108").Value = 0.6819444444416667
$ws.Range("AC109:AC116").Value = 0.6541666666666667
$ws.Range("AC117:AC146").Value = 0.4166666666666667
$ws.Range("AC147:AC176").Value = 0.4916666666666667
$ws.Range("AC177:AC221").Value = 0.4083333333333334

# --- 3. Append twelve new daily rows (222-233) ---
$newDates = @(
    "9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020",
    "10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020",
    "10/10/2020","10/11/2020"
)

# Same closure pattern (columns B..AB) as the prior last row (221), and the
# corresponding recomputed LockdownEffectiveness (column AC).
$pattern = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,0,0,0,0,0)
$acValue = 0.4083333333333334

$row = 222
foreach ($d in $newDates) {
    # Write the date label as literal text (not an auto-converted date
    # serial) by entering it as a string-literal formula, then flattening
    # the formula down to its plain value.
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = '="' + $d + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    for ($i = 0; $i -lt $pattern.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $pattern[$i]
    }
    $ws.Cells.Item($row, 29).Value = $acValue

    $row++
}

# Match the styling (bold / centered / bordered) already used by the other
# column-A date labels.
$ws.Range("A221").Copy()
$ws.Range("A222:A233").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
